$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 99993
$ws.Cells.Item(28, 3).Value = "Niet te lokaliseren"
$ws.Cells.Item(28, 4).Value = "Niet te lokaliseren"

$ws.Cells.Item(28, 1).NumberFormat = $ws.Cells.Item(27, 1).NumberFormat
$ws.Cells.Item(28, 2).NumberFormat = $ws.Cells.Item(27, 2).NumberFormat
